# Update NATMI LR-pair data (Fgf22-Fgfr2) with newly recomputed TPM-derived values.
# The previous "Resolving-Mac" sending-cluster block (rows 11-13) is dropped entirely,
# and the "Resolving-Mac" cluster used elsewhere in the table is renamed to
# "Inflammatory-Mac", while every numeric column is refreshed with new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing "Resolving-Mac" rows (11-13) and shift remaining rows up,
# shrinking the sheet range from A1:T13 to A1:T10.
$ws.Range("A11:T13").Delete(-4162) | Out-Null

# Row 2: FAPs (sending) -> ECs (target)
$ws.Range("A2").Value2 = "FAPs"
$ws.Range("B2").Value2 = "Fgf22"
$ws.Range("C2").Value2 = "Fgfr2"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 0.6666666666666666
$ws.Range("G2").Value2 = 0.2909523333333333
$ws.Range("H2").Value2 = 0.872857
$ws.Range("I2").Value2 = 0.68337437121998
$ws.Range("J2").Value2 = 0.68337437121998
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 0.106124
$ws.Range("N2").Value2 = 0.318372
$ws.Range("O2").Value2 = 0.08094716512538251
$ws.Range("P2").Value2 = 0.08094716512538253
$ws.Range("Q2").Value2 = 0.03087702542266666
$ws.Range("R2").Value2 = 0.277893228804
$ws.Range("S2").Value2 = 0.05531721806959817
$ws.Range("T2").Value2 = 0.05531721806959818

# Row 3: FAPs (sending) -> FAPs (target)
$ws.Range("A3").Value2 = "FAPs"
$ws.Range("B3").Value2 = "Fgf22"
$ws.Range("C3").Value2 = "Fgfr2"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 2
$ws.Range("F3").Value2 = 0.6666666666666666
$ws.Range("G3").Value2 = 0.2909523333333333
$ws.Range("H3").Value2 = 0.872857
$ws.Range("I3").Value2 = 0.68337437121998
$ws.Range("J3").Value2 = 0.68337437121998
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 1.092289666666667
$ws.Range("N3").Value2 = 3.276869
$ws.Range("O3").Value2 = 0.8331551016962769
$ws.Range("P3").Value2 = 0.833155101696277
$ws.Range("Q3").Value2 = 0.3178042271925555
$ws.Range("R3").Value2 = 2.860238044733
$ws.Range("S3").Value2 = 0.5693568437504117
$ws.Range("T3").Value2 = 0.5693568437504118

# Row 4: FAPs (sending) -> MuSCs (target)
$ws.Range("A4").Value2 = "FAPs"
$ws.Range("B4").Value2 = "Fgf22"
$ws.Range("C4").Value2 = "Fgfr2"
$ws.Range("D4").Value2 = "MuSCs"
$ws.Range("E4").Value2 = 2
$ws.Range("F4").Value2 = 0.6666666666666666
$ws.Range("G4").Value2 = 0.2909523333333333
$ws.Range("H4").Value2 = 0.872857
$ws.Range("I4").Value2 = 0.68337437121998
$ws.Range("J4").Value2 = 0.68337437121998
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 0.1126143333333333
$ws.Range("N4").Value2 = 0.337843
$ws.Range("O4").Value2 = 0.08589773317834044
$ws.Range("P4").Value2 = 0.08589773317834046
$ws.Range("Q4").Value2 = 0.03276540305011111
$ws.Range("R4").Value2 = 0.294888627451
$ws.Range("S4").Value2 = 0.05870030939997001
$ws.Range("T4").Value2 = 0.05870030939997003

# Row 5: Inflammatory-Mac (sending) -> ECs (target)
$ws.Range("A5").Value2 = "Inflammatory-Mac"
$ws.Range("B5").Value2 = "Fgf22"
$ws.Range("C5").Value2 = "Fgfr2"
$ws.Range("D5").Value2 = "ECs"
$ws.Range("E5").Value2 = 1
$ws.Range("F5").Value2 = 0.3333333333333333
$ws.Range("G5").Value2 = 0.07573966666666666
$ws.Range("H5").Value2 = 0.227219
$ws.Range("I5").Value2 = 0.1778935624669707
$ws.Range("J5").Value2 = 0.1778935624669707
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 0.106124
$ws.Range("N5").Value2 = 0.318372
$ws.Range("O5").Value2 = 0.08094716512538251
$ws.Range("P5").Value2 = 0.08094716512538253
$ws.Range("Q5").Value2 = 0.008037796385333333
$ws.Range("R5").Value2 = 0.072340167468
$ws.Range("S5").Value2 = 0.01439997957575643
$ws.Range("T5").Value2 = 0.01439997957575643

# Row 6: Inflammatory-Mac (sending) -> FAPs (target)
$ws.Range("A6").Value2 = "Inflammatory-Mac"
$ws.Range("B6").Value2 = "Fgf22"
$ws.Range("C6").Value2 = "Fgfr2"
$ws.Range("D6").Value2 = "FAPs"
$ws.Range("E6").Value2 = 1
$ws.Range("F6").Value2 = 0.3333333333333333
$ws.Range("G6").Value2 = 0.07573966666666666
$ws.Range("H6").Value2 = 0.227219
$ws.Range("I6").Value2 = 0.1778935624669707
$ws.Range("J6").Value2 = 0.1778935624669707
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 1.092289666666667
$ws.Range("N6").Value2 = 3.276869
$ws.Range("O6").Value2 = 0.8331551016962769
$ws.Range("P6").Value2 = 0.833155101696277
$ws.Range("Q6").Value2 = 0.08272965525677776
$ws.Range("R6").Value2 = 0.744566897311
$ws.Range("S6").Value2 = 0.148212929128282
$ws.Range("T6").Value2 = 0.148212929128282

# Row 7: Inflammatory-Mac (sending) -> MuSCs (target)
$ws.Range("A7").Value2 = "Inflammatory-Mac"
$ws.Range("B7").Value2 = "Fgf22"
$ws.Range("C7").Value2 = "Fgfr2"
$ws.Range("D7").Value2 = "MuSCs"
$ws.Range("E7").Value2 = 1
$ws.Range("F7").Value2 = 0.3333333333333333
$ws.Range("G7").Value2 = 0.07573966666666666
$ws.Range("H7").Value2 = 0.227219
$ws.Range("I7").Value2 = 0.1778935624669707
$ws.Range("J7").Value2 = 0.1778935624669707
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 0.1126143333333333
$ws.Range("N7").Value2 = 0.337843
$ws.Range("O7").Value2 = 0.08589773317834044
$ws.Range("P7").Value2 = 0.08589773317834046
$ws.Range("Q7").Value2 = 0.008529372068555556
$ws.Range("R7").Value2 = 0.076764348617
$ws.Range("S7").Value2 = 0.01528065376293229
$ws.Range("T7").Value2 = 0.01528065376293229

# Row 8: MuSCs (sending) -> ECs (target)
$ws.Range("A8").Value2 = "MuSCs"
$ws.Range("B8").Value2 = "Fgf22"
$ws.Range("C8").Value2 = "Fgfr2"
$ws.Range("D8").Value2 = "ECs"
$ws.Range("E8").Value2 = 1
$ws.Range("F8").Value2 = 0.3333333333333333
$ws.Range("G8").Value2 = 0.05906633333333333
$ws.Range("H8").Value2 = 0.177199
$ws.Range("I8").Value2 = 0.1387320663130493
$ws.Range("J8").Value2 = 0.1387320663130493
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 0.106124
$ws.Range("N8").Value2 = 0.318372
$ws.Range("O8").Value2 = 0.08094716512538251
$ws.Range("P8").Value2 = 0.08094716512538253
$ws.Range("Q8").Value2 = 0.006268355558666667
$ws.Range("R8").Value2 = 0.056415200028
$ws.Range("S8").Value2 = 0.01122996748002792
$ws.Range("T8").Value2 = 0.01122996748002792

# Row 9: MuSCs (sending) -> FAPs (target)
$ws.Range("A9").Value2 = "MuSCs"
$ws.Range("B9").Value2 = "Fgf22"
$ws.Range("C9").Value2 = "Fgfr2"
$ws.Range("D9").Value2 = "FAPs"
$ws.Range("E9").Value2 = 1
$ws.Range("F9").Value2 = 0.3333333333333333
$ws.Range("G9").Value2 = 0.05906633333333333
$ws.Range("H9").Value2 = 0.177199
$ws.Range("I9").Value2 = 0.1387320663130493
$ws.Range("J9").Value2 = 0.1387320663130493
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 1.092289666666667
$ws.Range("N9").Value2 = 3.276869
$ws.Range("O9").Value2 = 0.8331551016962769
$ws.Range("P9").Value2 = 0.833155101696277
$ws.Range("Q9").Value2 = 0.06451754554788888
$ws.Range("R9").Value2 = 0.580657909931
$ws.Range("S9").Value2 = 0.1155853288175832
$ws.Range("T9").Value2 = 0.1155853288175832

# Row 10: MuSCs (sending) -> MuSCs (target)
$ws.Range("A10").Value2 = "MuSCs"
$ws.Range("B10").Value2 = "Fgf22"
$ws.Range("C10").Value2 = "Fgfr2"
$ws.Range("D10").Value2 = "MuSCs"
$ws.Range("E10").Value2 = 1
$ws.Range("F10").Value2 = 0.3333333333333333
$ws.Range("G10").Value2 = 0.05906633333333333
$ws.Range("H10").Value2 = 0.177199
$ws.Range("I10").Value2 = 0.1387320663130493
$ws.Range("J10").Value2 = 0.1387320663130493
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 0.1126143333333333
$ws.Range("N10").Value2 = 0.337843
$ws.Range("O10").Value2 = 0.08589773317834044
$ws.Range("P10").Value2 = 0.08589773317834046
$ws.Range("Q10").Value2 = 0.006651715750777778
$ws.Range("R10").Value2 = 0.059865441757
$ws.Range("S10").Value2 = 0.01191677001543814
$ws.Range("T10").Value2 = 0.01191677001543814
